$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ESO load data for rows 5-9 (columns D:AA), plus a new appended
# row 10 for 2026-02-13. Columns A:C (Year/Month/Day) are unchanged for
# the existing rows, so only D:AA are rewritten there.

$updatedLoadRows = @{
    5 = @(4703, 4461, 4349, 4329, 4380, 4395, 4480, 4665, 4950, 5225, 5392, 5467, 5508, 5377, 5377, 5400, 5504, 5699, 5903, 5798, 5585, 5345, 5170, 4881)
    6 = @(4954, 4717, 4610, 4617, 4622, 4721, 5017, 5437, 5872, 6121, 6188, 6171, 6126, 6097, 6047, 5945, 5866, 5982, 6324, 6311, 6123, 5850, 5581, 5261)
    7 = @(5257, 5034, 4944, 4937, 4936, 5066, 5469, 5930, 6311, 6439, 6367, 6265, 6168, 6116, 6046, 6059, 6104, 6273, 6567, 6540, 6337, 6058, 5836, 5518)
    8 = @(5231, 5008, 4918, 4911, 4910, 5040, 5443, 5903, 6285, 6412, 6342, 6242, 6147, 6095, 6027, 6040, 6084, 6251, 6541, 6514, 6311, 6032, 5810, 5491)
    9 = @(5036, 4812, 4721, 4714, 4713, 4842, 5245, 5705, 6086, 6214, 6154, 6069, 5987, 5943, 5885, 5896, 5937, 6087, 6349, 6322, 6117, 5834, 5610, 5288)
}

foreach ($rowNum in $updatedLoadRows.Keys) {
    $values = $updatedLoadRows[$rowNum]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($rowNum, 4 + $i).Value = $values[$i]
    }
}

$newRow10 = @(2026, 2, 13, 4852, 4627, 4536, 4529, 4528, 4657, 5059, 5519, 5900, 6027, 5977, 5905, 5837, 5800, 5751, 5760, 5789, 5898, 6089, 6062, 5861, 5584, 5364, 5048)
for ($i = 0; $i -lt $newRow10.Length; $i++) {
    $ws.Cells.Item(10, 1 + $i).Value = $newRow10[$i]
}
